# Move the "X" mark from the "Nee" (No) column to the "Ja" (Yes) column
# for the row "Je hebt opdrachten in canvas gemaakt ter voorbereiding" in
# the first table of the document.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$targetText = "Je hebt opdrachten in canvas gemaakt ter voorbereiding"

function TrimMarks($s) {
    return $s.TrimEnd([char]7).TrimEnd([char]13)
}

$targetRow = $null
$rowCount = $t.Rows.Count
for ($i = 1; $i -le $rowCount; $i++) {
    $row = $t.Rows.Item($i)
    if ($row.Cells.Count -ge 3) {
        $c1Text = TrimMarks $row.Cells.Item(1).Range.Text
        if ($c1Text -eq $targetText) {
            $targetRow = $row
            break
        }
    }
}

if ($targetRow -ne $null) {
    $jaCell = $targetRow.Cells.Item(2)
    $neeCell = $targetRow.Cells.Item(3)

    # Add "X" to the "Ja" cell, right after its existing content.
    $jaRange = $jaCell.Range
    $jaClean = TrimMarks $jaRange.Text
    $jaRange.Text = $jaClean + "X"

    # Remove all visible text from the "Nee" cell, leaving a bare empty
    # paragraph (no runs at all) behind - mirrors it being unchecked.
    $guard = 0
    while ((TrimMarks $neeCell.Range.Text).Length -gt 0 -and $guard -lt 50) {
        $r = $neeCell.Range
        $one = $d.Range($r.Start, $r.Start + 1)
        $one.Text = ""
        $guard = $guard + 1
    }
}
